$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed by Excel as a number;
# force them to remain plain text (matching the source inlineStr cells) by temporarily
# switching to a text number format, then restoring the original style afterwards.
function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '69.612.65'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.503.94'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '574.47'
$ws.Range('E5').Value = '  -0.57%  '
Set-TextValue $ws.Range('D6') '166.35'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('E7').Value = '  -0.05%  '
Set-TextValue $ws.Range('D8') '0.513'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').Value = '2.502.40'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  +3.15%  '
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '2.958.29'
$ws.Range('D15').Value = '69.480.89'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('E16').Value = '  +0.46%  '
Set-TextValue $ws.Range('D17') '24.70'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').Value = '2.486.51'
$ws.Range('E18').Value = '  -1.47%  '
Set-TextValue $ws.Range('D19') '11.21'
$ws.Range('E19').Value = '  -1.08%  '
Set-TextValue $ws.Range('D20') '7.51'
$ws.Range('E20').Value = '  -3.54%  '
Set-TextValue $ws.Range('D21') '348.88'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('E22').Value = '  -1.09%  '
Set-TextValue $ws.Range('D23') '1.94'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('E24').Value = '  +0.03%  '
Set-TextValue $ws.Range('D25') '70.69'
$ws.Range('E25').Value = '  +2.16%  '
Set-TextValue $ws.Range('D26') '3.93'
$ws.Range('E26').Value = '  -2.14%  '
Set-TextValue $ws.Range('D27') '8.73'
$ws.Range('E27').Value = '  -3.28%  '
$ws.Range('D28').Value = '2.630.07'
$ws.Range('E28').Value = '  -0.69%  '
Set-TextValue $ws.Range('D29') '0.997'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('E30').Value = '  -2.31%  '
Set-TextValue $ws.Range('D31') '7.86'
$ws.Range('E31').Value = '  -0.65%  '
Set-TextValue $ws.Range('D32') '457.52'
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('E33').Value = '  -5.90%  '
$ws.Range('E34').Value = '  -1.62%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D36') '0.116'
$ws.Range('E36').Value = '  -3.85%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D37') '157.21'
$ws.Range('E37').Value = '  +2.75%  '
Set-TextValue $ws.Range('D38') '19.04'
$ws.Range('E38').Value = '  +0.07%  '
Set-TextValue $ws.Range('D39') '18.35'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('E41').Value = '  -1.11%  '
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('E43').Value = '  -0.32%  '
Set-TextValue $ws.Range('D44') '38.15'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('E45').Value = '  -5.37%  '
$ws.Range('E46').Value = '  -7.98%  '
Set-TextValue $ws.Range('D47') '141.09'
$ws.Range('E47').Value = '  -1.62%  '
Set-TextValue $ws.Range('D48') '3.49'
$ws.Range('E48').Value = '  -0.59%  '
Set-TextValue $ws.Range('D49') '0.519'
$ws.Range('E49').Value = '  -2.69%  '
Set-TextValue $ws.Range('D50') '0.0732'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('E51').Value = '  -0.90%  '
